$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "예윤 병국"
$ws.Range("C1").Value = "예윤 병국"
$ws.Range("B3").Value = "태훈 현빈"
$ws.Range("B4").Value = "태훈 태훈"
